$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The clan log was refreshed: two new recruits appended (ergo / Ergo460),
# trophies/donations re-synced for the existing roster, and the sheet
# re-sorted by Trophies (desc) same as the generator always does.
# First, extend formatting for the two brand-new rows (45 data rows -> 47)
# by cloning the zebra-stripe style from the last two existing rows.
$ws.Range("A44:I44").Copy()
$ws.Range("A46:I46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A45:I45").Copy()
$ws.Range("A47:I47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Full rewrite of data rows (A2:I47) to match final clan snapshot ---
$ws.Cells.Item(2,1).Value = 1.0
$ws.Cells.Item(2,2).Value = "Markś Village"
$ws.Cells.Item(2,3).Value = "#LRUGY0PQ"
$ws.Cells.Item(2,4).Value = 178.0
$ws.Cells.Item(2,5).Value = 5138.0
$ws.Cells.Item(2,6).Value = "admin"
$ws.Cells.Item(2,7).Value = 1816.0
$ws.Cells.Item(2,8).Value = 1461.0
$ws.Cells.Item(2,9).Value = 1.24

$ws.Cells.Item(3,1).Value = 2.0
$ws.Cells.Item(3,2).Value = "Titi"
$ws.Cells.Item(3,3).Value = "#2JP0VYLL"
$ws.Cells.Item(3,4).Value = 198.0
$ws.Cells.Item(3,5).Value = 5070.0
$ws.Cells.Item(3,6).Value = "coLeader"
$ws.Cells.Item(3,7).Value = 1112.0
$ws.Cells.Item(3,8).Value = 441.0
$ws.Cells.Item(3,9).Value = 2.52

$ws.Cells.Item(4,1).Value = 3.0
$ws.Cells.Item(4,2).Value = "Sem"
$ws.Cells.Item(4,3).Value = "#LVQJUYC"
$ws.Cells.Item(4,4).Value = 211.0
$ws.Cells.Item(4,5).Value = 5022.0
$ws.Cells.Item(4,6).Value = "leader"
$ws.Cells.Item(4,7).Value = 487.0
$ws.Cells.Item(4,8).Value = 370.0
$ws.Cells.Item(4,9).Value = 1.32

$ws.Cells.Item(5,1).Value = 4.0
$ws.Cells.Item(5,2).Value = "bastaard"
$ws.Cells.Item(5,3).Value = "#2JGULJY8"
$ws.Cells.Item(5,4).Value = 210.0
$ws.Cells.Item(5,5).Value = 4991.0
$ws.Cells.Item(5,6).Value = "coLeader"
$ws.Cells.Item(5,7).Value = 733.0
$ws.Cells.Item(5,8).Value = 695.0
$ws.Cells.Item(5,9).Value = 1.05

$ws.Cells.Item(6,1).Value = 5.0
$ws.Cells.Item(6,2).Value = "clan de renzo"
$ws.Cells.Item(6,3).Value = "#J298GYQC"
$ws.Cells.Item(6,4).Value = 180.0
$ws.Cells.Item(6,5).Value = 4985.0
$ws.Cells.Item(6,6).Value = "coLeader"
$ws.Cells.Item(6,7).Value = 324.0
$ws.Cells.Item(6,8).Value = 512.0
$ws.Cells.Item(6,9).Value = 0.63

$ws.Cells.Item(7,1).Value = 6.0
$ws.Cells.Item(7,2).Value = "ed"
$ws.Cells.Item(7,3).Value = "#QR2VVJJG"
$ws.Cells.Item(7,4).Value = 189.0
$ws.Cells.Item(7,5).Value = 4978.0
$ws.Cells.Item(7,6).Value = "coLeader"
$ws.Cells.Item(7,7).Value = 84.0
$ws.Cells.Item(7,8).Value = 558.0
$ws.Cells.Item(7,9).Value = 0.15

$ws.Cells.Item(8,1).Value = 7.0
$ws.Cells.Item(8,2).Value = "Bam"
$ws.Cells.Item(8,3).Value = "#YRJ8J8VG"
$ws.Cells.Item(8,4).Value = 179.0
$ws.Cells.Item(8,5).Value = 4967.0
$ws.Cells.Item(8,6).Value = "coLeader"
$ws.Cells.Item(8,7).Value = 230.0
$ws.Cells.Item(8,8).Value = 97.0
$ws.Cells.Item(8,9).Value = 2.37

$ws.Cells.Item(9,1).Value = 8.0
$ws.Cells.Item(9,2).Value = "billie jean"
$ws.Cells.Item(9,3).Value = "#QLV8Q0C0"
$ws.Cells.Item(9,4).Value = 172.0
$ws.Cells.Item(9,5).Value = 4924.0
$ws.Cells.Item(9,6).Value = "coLeader"
$ws.Cells.Item(9,7).Value = 646.0
$ws.Cells.Item(9,8).Value = 479.0
$ws.Cells.Item(9,9).Value = 1.35

$ws.Cells.Item(10,1).Value = 9.0
$ws.Cells.Item(10,2).Value = "Beertjuh#"
$ws.Cells.Item(10,3).Value = "#V8QY2UUJ"
$ws.Cells.Item(10,4).Value = 179.0
$ws.Cells.Item(10,5).Value = 4859.0
$ws.Cells.Item(10,6).Value = "coLeader"
$ws.Cells.Item(10,7).Value = 576.0
$ws.Cells.Item(10,8).Value = 364.0
$ws.Cells.Item(10,9).Value = 1.58

$ws.Cells.Item(11,1).Value = 10.0
$ws.Cells.Item(11,2).Value = "GJS"
$ws.Cells.Item(11,3).Value = "#8YJRVRLC"
$ws.Cells.Item(11,4).Value = 162.0
$ws.Cells.Item(11,5).Value = 4853.0
$ws.Cells.Item(11,6).Value = "admin"
$ws.Cells.Item(11,7).Value = 182.0
$ws.Cells.Item(11,8).Value = 331.0
$ws.Cells.Item(11,9).Value = 0.55

$ws.Cells.Item(12,1).Value = 11.0
$ws.Cells.Item(12,2).Value = "shag tand"
$ws.Cells.Item(12,3).Value = "#90VUQ988"
$ws.Cells.Item(12,4).Value = 184.0
$ws.Cells.Item(12,5).Value = 4847.0
$ws.Cells.Item(12,6).Value = "coLeader"
$ws.Cells.Item(12,7).Value = 277.0
$ws.Cells.Item(12,8).Value = 77.0
$ws.Cells.Item(12,9).Value = 3.6

$ws.Cells.Item(13,1).Value = 12.0
$ws.Cells.Item(13,2).Value = "tomtoch"
$ws.Cells.Item(13,3).Value = "#8299JUJY"
$ws.Cells.Item(13,4).Value = 180.0
$ws.Cells.Item(13,5).Value = 4841.0
$ws.Cells.Item(13,6).Value = "coLeader"
$ws.Cells.Item(13,7).Value = 225.0
$ws.Cells.Item(13,8).Value = 333.0
$ws.Cells.Item(13,9).Value = 0.68

$ws.Cells.Item(14,1).Value = 13.0
$ws.Cells.Item(14,2).Value = "wht^"
$ws.Cells.Item(14,3).Value = "#QU0VGY9"
$ws.Cells.Item(14,4).Value = 175.0
$ws.Cells.Item(14,5).Value = 4831.0
$ws.Cells.Item(14,6).Value = "coLeader"
$ws.Cells.Item(14,7).Value = 340.0
$ws.Cells.Item(14,8).Value = 292.0
$ws.Cells.Item(14,9).Value = 1.16

$ws.Cells.Item(15,1).Value = 14.0
$ws.Cells.Item(15,2).Value = "leeuw"
$ws.Cells.Item(15,3).Value = "#PP0C92QP"
$ws.Cells.Item(15,4).Value = 175.0
$ws.Cells.Item(15,5).Value = 4790.0
$ws.Cells.Item(15,6).Value = "coLeader"
$ws.Cells.Item(15,7).Value = 469.0
$ws.Cells.Item(15,8).Value = 498.0
$ws.Cells.Item(15,9).Value = 0.94

$ws.Cells.Item(16,1).Value = 15.0
$ws.Cells.Item(16,2).Value = "* haakie *"
$ws.Cells.Item(16,3).Value = "#2CUU0VJG"
$ws.Cells.Item(16,4).Value = 170.0
$ws.Cells.Item(16,5).Value = 4705.0
$ws.Cells.Item(16,6).Value = "admin"
$ws.Cells.Item(16,7).Value = 979.0
$ws.Cells.Item(16,8).Value = 406.0
$ws.Cells.Item(16,9).Value = 2.41

$ws.Cells.Item(17,1).Value = 16.0
$ws.Cells.Item(17,2).Value = "Z!P"
$ws.Cells.Item(17,3).Value = "#YVYU9Y8R"
$ws.Cells.Item(17,4).Value = 179.0
$ws.Cells.Item(17,5).Value = 4702.0
$ws.Cells.Item(17,6).Value = "coLeader"
$ws.Cells.Item(17,7).Value = 157.0
$ws.Cells.Item(17,8).Value = 140.0
$ws.Cells.Item(17,9).Value = 1.12

$ws.Cells.Item(18,1).Value = 17.0
$ws.Cells.Item(18,2).Value = "Ronald O"
$ws.Cells.Item(18,3).Value = "#99UJ999G"
$ws.Cells.Item(18,4).Value = 163.0
$ws.Cells.Item(18,5).Value = 4632.0
$ws.Cells.Item(18,6).Value = "coLeader"
$ws.Cells.Item(18,7).Value = 368.0
$ws.Cells.Item(18,8).Value = 646.0
$ws.Cells.Item(18,9).Value = 0.57

$ws.Cells.Item(19,1).Value = 18.0
$ws.Cells.Item(19,2).Value = "xavier"
$ws.Cells.Item(19,3).Value = "#8P88LUCQ"
$ws.Cells.Item(19,4).Value = 182.0
$ws.Cells.Item(19,5).Value = 4609.0
$ws.Cells.Item(19,6).Value = "coLeader"
$ws.Cells.Item(19,7).Value = 800.0
$ws.Cells.Item(19,8).Value = 860.0
$ws.Cells.Item(19,9).Value = 0.93

$ws.Cells.Item(20,1).Value = 19.0
$ws.Cells.Item(20,2).Value = "THICK GIRLS"
$ws.Cells.Item(20,3).Value = "#8U09PR0V"
$ws.Cells.Item(20,4).Value = 184.0
$ws.Cells.Item(20,5).Value = 4584.0
$ws.Cells.Item(20,6).Value = "coLeader"
$ws.Cells.Item(20,7).Value = 391.0
$ws.Cells.Item(20,8).Value = 148.0
$ws.Cells.Item(20,9).Value = 2.64

$ws.Cells.Item(21,1).Value = 20.0
$ws.Cells.Item(21,2).Value = "rommensje"
$ws.Cells.Item(21,3).Value = "#9Y0URQUL"
$ws.Cells.Item(21,4).Value = 192.0
$ws.Cells.Item(21,5).Value = 4532.0
$ws.Cells.Item(21,6).Value = "admin"
$ws.Cells.Item(21,7).Value = 226.0
$ws.Cells.Item(21,8).Value = 220.0
$ws.Cells.Item(21,9).Value = 1.03

$ws.Cells.Item(22,1).Value = 21.0
$ws.Cells.Item(22,2).Value = "duke en wesj"
$ws.Cells.Item(22,3).Value = "#P0PV2CJY"
$ws.Cells.Item(22,4).Value = 211.0
$ws.Cells.Item(22,5).Value = 4493.0
$ws.Cells.Item(22,6).Value = "coLeader"
$ws.Cells.Item(22,7).Value = 416.0
$ws.Cells.Item(22,8).Value = 1329.0
$ws.Cells.Item(22,9).Value = 0.31

$ws.Cells.Item(23,1).Value = 22.0
$ws.Cells.Item(23,2).Value = "Bastos"
$ws.Cells.Item(23,3).Value = "#8RP8QV8V"
$ws.Cells.Item(23,4).Value = 172.0
$ws.Cells.Item(23,5).Value = 4490.0
$ws.Cells.Item(23,6).Value = "member"
$ws.Cells.Item(23,7).Value = 0.0
$ws.Cells.Item(23,8).Value = 0.0
$ws.Cells.Item(23,9).Value = 0.0

$ws.Cells.Item(24,1).Value = 23.0
$ws.Cells.Item(24,2).Value = "Wolverine"
$ws.Cells.Item(24,3).Value = "#828YUV9J"
$ws.Cells.Item(24,4).Value = 161.0
$ws.Cells.Item(24,5).Value = 4440.0
$ws.Cells.Item(24,6).Value = "member"
$ws.Cells.Item(24,7).Value = 30.0
$ws.Cells.Item(24,8).Value = 57.0
$ws.Cells.Item(24,9).Value = 0.53

$ws.Cells.Item(25,1).Value = 24.0
$ws.Cells.Item(25,2).Value = "Shyngalicious"
$ws.Cells.Item(25,3).Value = "#RLCLPJ"
$ws.Cells.Item(25,4).Value = 149.0
$ws.Cells.Item(25,5).Value = 4437.0
$ws.Cells.Item(25,6).Value = "coLeader"
$ws.Cells.Item(25,7).Value = 267.0
$ws.Cells.Item(25,8).Value = 181.0
$ws.Cells.Item(25,9).Value = 1.48

$ws.Cells.Item(26,1).Value = 25.0
$ws.Cells.Item(26,2).Value = "elandro"
$ws.Cells.Item(26,3).Value = "#22GU992L"
$ws.Cells.Item(26,4).Value = 175.0
$ws.Cells.Item(26,5).Value = 4399.0
$ws.Cells.Item(26,6).Value = "coLeader"
$ws.Cells.Item(26,7).Value = 492.0
$ws.Cells.Item(26,8).Value = 706.0
$ws.Cells.Item(26,9).Value = 0.7

$ws.Cells.Item(27,1).Value = 26.0
$ws.Cells.Item(27,2).Value = "BlackWing"
$ws.Cells.Item(27,3).Value = "#PJ8CG2J9"
$ws.Cells.Item(27,4).Value = 160.0
$ws.Cells.Item(27,5).Value = 4269.0
$ws.Cells.Item(27,6).Value = "coLeader"
$ws.Cells.Item(27,7).Value = 72.0
$ws.Cells.Item(27,8).Value = 227.0
$ws.Cells.Item(27,9).Value = 0.32

$ws.Cells.Item(28,1).Value = 27.0
$ws.Cells.Item(28,2).Value = "peter"
$ws.Cells.Item(28,3).Value = "#8LV09JG9"
$ws.Cells.Item(28,4).Value = 160.0
$ws.Cells.Item(28,5).Value = 4262.0
$ws.Cells.Item(28,6).Value = "coLeader"
$ws.Cells.Item(28,7).Value = 180.0
$ws.Cells.Item(28,8).Value = 313.0
$ws.Cells.Item(28,9).Value = 0.58

$ws.Cells.Item(29,1).Value = 28.0
$ws.Cells.Item(29,2).Value = "pamuk39"
$ws.Cells.Item(29,3).Value = "#RV8JG08P"
$ws.Cells.Item(29,4).Value = 195.0
$ws.Cells.Item(29,5).Value = 4259.0
$ws.Cells.Item(29,6).Value = "coLeader"
$ws.Cells.Item(29,7).Value = 951.0
$ws.Cells.Item(29,8).Value = 852.0
$ws.Cells.Item(29,9).Value = 1.12

$ws.Cells.Item(30,1).Value = 29.0
$ws.Cells.Item(30,2).Value = "(j)de tik(j)"
$ws.Cells.Item(30,3).Value = "#GYVQ0Y8R"
$ws.Cells.Item(30,4).Value = 178.0
$ws.Cells.Item(30,5).Value = 4233.0
$ws.Cells.Item(30,6).Value = "coLeader"
$ws.Cells.Item(30,7).Value = 126.0
$ws.Cells.Item(30,8).Value = 222.0
$ws.Cells.Item(30,9).Value = 0.57

$ws.Cells.Item(31,1).Value = 30.0
$ws.Cells.Item(31,2).Value = "mauzer99"
$ws.Cells.Item(31,3).Value = "#P8VQ09QQ"
$ws.Cells.Item(31,4).Value = 172.0
$ws.Cells.Item(31,5).Value = 4125.0
$ws.Cells.Item(31,6).Value = "member"
$ws.Cells.Item(31,7).Value = 0.0
$ws.Cells.Item(31,8).Value = 112.0
$ws.Cells.Item(31,9).Value = 0.0

$ws.Cells.Item(32,1).Value = 31.0
$ws.Cells.Item(32,2).Value = "kevintjuh93"
$ws.Cells.Item(32,3).Value = "#YG0URYGQ"
$ws.Cells.Item(32,4).Value = 146.0
$ws.Cells.Item(32,5).Value = 4121.0
$ws.Cells.Item(32,6).Value = "member"
$ws.Cells.Item(32,7).Value = 288.0
$ws.Cells.Item(32,8).Value = 329.0
$ws.Cells.Item(32,9).Value = 0.88

$ws.Cells.Item(33,1).Value = 32.0
$ws.Cells.Item(33,2).Value = "axes"
$ws.Cells.Item(33,3).Value = "#2JVRYC22"
$ws.Cells.Item(33,4).Value = 168.0
$ws.Cells.Item(33,5).Value = 3842.0
$ws.Cells.Item(33,6).Value = "admin"
$ws.Cells.Item(33,7).Value = 177.0
$ws.Cells.Item(33,8).Value = 430.0
$ws.Cells.Item(33,9).Value = 0.41

$ws.Cells.Item(34,1).Value = 33.0
$ws.Cells.Item(34,2).Value = "wiski"
$ws.Cells.Item(34,3).Value = "#P8LL80GV"
$ws.Cells.Item(34,4).Value = 156.0
$ws.Cells.Item(34,5).Value = 3803.0
$ws.Cells.Item(34,6).Value = "admin"
$ws.Cells.Item(34,7).Value = 283.0
$ws.Cells.Item(34,8).Value = 140.0
$ws.Cells.Item(34,9).Value = 2.02

$ws.Cells.Item(35,1).Value = 34.0
$ws.Cells.Item(35,2).Value = "PeterClash"
$ws.Cells.Item(35,3).Value = "#2VGG9R288"
$ws.Cells.Item(35,4).Value = 171.0
$ws.Cells.Item(35,5).Value = 3726.0
$ws.Cells.Item(35,6).Value = "coLeader"
$ws.Cells.Item(35,7).Value = 847.0
$ws.Cells.Item(35,8).Value = 573.0
$ws.Cells.Item(35,9).Value = 1.48

$ws.Cells.Item(36,1).Value = 35.0
$ws.Cells.Item(36,2).Value = "Rolex"
$ws.Cells.Item(36,3).Value = "#QJLQY0VY"
$ws.Cells.Item(36,4).Value = 159.0
$ws.Cells.Item(36,5).Value = 3598.0
$ws.Cells.Item(36,6).Value = "member"
$ws.Cells.Item(36,7).Value = 0.0
$ws.Cells.Item(36,8).Value = 0.0
$ws.Cells.Item(36,9).Value = 0.0

$ws.Cells.Item(37,1).Value = 36.0
$ws.Cells.Item(37,2).Value = "Lol-with-pohl"
$ws.Cells.Item(37,3).Value = "#G2UGVYCJ"
$ws.Cells.Item(37,4).Value = 166.0
$ws.Cells.Item(37,5).Value = 3435.0
$ws.Cells.Item(37,6).Value = "admin"
$ws.Cells.Item(37,7).Value = 504.0
$ws.Cells.Item(37,8).Value = 413.0
$ws.Cells.Item(37,9).Value = 1.22

$ws.Cells.Item(38,1).Value = 37.0
$ws.Cells.Item(38,2).Value = "•#FRANK#•"
$ws.Cells.Item(38,3).Value = "#28YP9RL0G"
$ws.Cells.Item(38,4).Value = 132.0
$ws.Cells.Item(38,5).Value = 3368.0
$ws.Cells.Item(38,6).Value = "member"
$ws.Cells.Item(38,7).Value = 0.0
$ws.Cells.Item(38,8).Value = 0.0
$ws.Cells.Item(38,9).Value = 0.0

$ws.Cells.Item(39,1).Value = 38.0
$ws.Cells.Item(39,2).Value = "KLEINE EDJE"
$ws.Cells.Item(39,3).Value = "#C0UC2QPJ"
$ws.Cells.Item(39,4).Value = 151.0
$ws.Cells.Item(39,5).Value = 3345.0
$ws.Cells.Item(39,6).Value = "coLeader"
$ws.Cells.Item(39,7).Value = 791.0
$ws.Cells.Item(39,8).Value = 30.0
$ws.Cells.Item(39,9).Value = 26.37

$ws.Cells.Item(40,1).Value = 39.0
$ws.Cells.Item(40,2).Value = "marco25"
$ws.Cells.Item(40,3).Value = "#9CGP2Y2G"
$ws.Cells.Item(40,4).Value = 214.0
$ws.Cells.Item(40,5).Value = 3309.0
$ws.Cells.Item(40,6).Value = "admin"
$ws.Cells.Item(40,7).Value = 872.0
$ws.Cells.Item(40,8).Value = 915.0
$ws.Cells.Item(40,9).Value = 0.95

$ws.Cells.Item(41,1).Value = 40.0
$ws.Cells.Item(41,2).Value = "SUPREMACYYYY"
$ws.Cells.Item(41,3).Value = "#8VP0J2QP"
$ws.Cells.Item(41,4).Value = 152.0
$ws.Cells.Item(41,5).Value = 3248.0
$ws.Cells.Item(41,6).Value = "member"
$ws.Cells.Item(41,7).Value = 302.0
$ws.Cells.Item(41,8).Value = 460.0
$ws.Cells.Item(41,9).Value = 0.66

$ws.Cells.Item(42,1).Value = 41.0
$ws.Cells.Item(42,2).Value = "ergo"
$ws.Cells.Item(42,3).Value = "#QR2LJQY"
$ws.Cells.Item(42,4).Value = 154.0
$ws.Cells.Item(42,5).Value = 3211.0
$ws.Cells.Item(42,6).Value = "coLeader"
$ws.Cells.Item(42,7).Value = 4.0
$ws.Cells.Item(42,8).Value = 35.0
$ws.Cells.Item(42,9).Value = 0.11

$ws.Cells.Item(43,1).Value = 42.0
$ws.Cells.Item(43,2).Value = "OhhJayky!"
$ws.Cells.Item(43,3).Value = "#R0Y2J8G2"
$ws.Cells.Item(43,4).Value = 123.0
$ws.Cells.Item(43,5).Value = 3199.0
$ws.Cells.Item(43,6).Value = "admin"
$ws.Cells.Item(43,7).Value = 69.0
$ws.Cells.Item(43,8).Value = 62.0
$ws.Cells.Item(43,9).Value = 1.11

$ws.Cells.Item(44,1).Value = 43.0
$ws.Cells.Item(44,2).Value = "jojo"
$ws.Cells.Item(44,3).Value = "#22G9G88CC"
$ws.Cells.Item(44,4).Value = 156.0
$ws.Cells.Item(44,5).Value = 3113.0
$ws.Cells.Item(44,6).Value = "admin"
$ws.Cells.Item(44,7).Value = 407.0
$ws.Cells.Item(44,8).Value = 988.0
$ws.Cells.Item(44,9).Value = 0.41

$ws.Cells.Item(45,1).Value = 44.0
$ws.Cells.Item(45,2).Value = "Ergo460"
$ws.Cells.Item(45,3).Value = "#98LPPUJQ"
$ws.Cells.Item(45,4).Value = 136.0
$ws.Cells.Item(45,5).Value = 2656.0
$ws.Cells.Item(45,6).Value = "admin"
$ws.Cells.Item(45,7).Value = 0.0
$ws.Cells.Item(45,8).Value = 3.0
$ws.Cells.Item(45,9).Value = 0.0

$ws.Cells.Item(46,1).Value = 45.0
$ws.Cells.Item(46,2).Value = "Duke en wes"
$ws.Cells.Item(46,3).Value = "#8G2090PC9"
$ws.Cells.Item(46,4).Value = 113.0
$ws.Cells.Item(46,5).Value = 2577.0
$ws.Cells.Item(46,6).Value = "coLeader"
$ws.Cells.Item(46,7).Value = 124.0
$ws.Cells.Item(46,8).Value = 184.0
$ws.Cells.Item(46,9).Value = 0.67

$ws.Cells.Item(47,1).Value = 46.0
$ws.Cells.Item(47,2).Value = "Raffie"
$ws.Cells.Item(47,3).Value = "#PY0QCYUPR"
$ws.Cells.Item(47,4).Value = 69.0
$ws.Cells.Item(47,5).Value = 1333.0
$ws.Cells.Item(47,6).Value = "member"
$ws.Cells.Item(47,7).Value = 0.0
$ws.Cells.Item(47,8).Value = 40.0
$ws.Cells.Item(47,9).Value = 0.0


# Update the printed footer timestamp to match the regenerated report.
$ps = $ws.PageSetup
$ps.LeftFooter = "Clanoverzicht"
$ps.RightFooter = "27/12/2017 16:05"
